$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------
# "Айди" (old/typo transliteration) -> "ИД" (new abbreviation) in the
# Bulgarian "SigiDoc ID" labels.
$ws.Range("B7").Value = "СигиДок ИД: 13"
$ws.Range("C7").Value = "СигиДок ИД: 11"

# "latitude" / "longitude" labels uppercased to match the other field
# labels (e.g. SEAL ID, TYPE, ...), and promoted to the section-header
# style used by the other header rows (copy format only, keep values).
$ws.Range("A100").Copy() | Out-Null
$ws.Range("A106").PasteSpecial(-4122) | Out-Null
$ws.Range("A106").Value = "LATITUDE"

$ws.Range("A101").Copy() | Out-Null
$ws.Range("A107").PasteSpecial(-4122) | Out-Null
$ws.Range("A107").Value = "LONGITUDE"

$excel.CutCopyMode = 0

# --- Row heights -------------------------------------------------------
# Rows 100/101 revert to the sheet's default row height.
$ws.Rows.Item(100).RowHeight = $ws.Rows.Item(1).RowHeight
$ws.Rows.Item(101).RowHeight = $ws.Rows.Item(1).RowHeight

# --- View state ----------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C7").Select()
